$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Kohler"
$ws.Range("B2").Value = "Alina"

$ws.Range("A3").Value = "Kohler"
$ws.Range("B3").Value = "Nina"

$ws.Range("A4").Value = "Matumona"
$ws.Range("B4").Value = "Noe"

$ws.Range("A5").Value = "Matumona"
$ws.Range("B5").Value = "Nina"

$ws.Range("A6").Value = "Sarman"
$ws.Range("B6").Value = "Dominik"

$ws.Range("A7").Value = "Zillig"
$ws.Range("B7").Value = "Nicolas"

$ws.Range("A8").Value = "asdf"
$ws.Range("B8").Value = "Marlene"
